# MassWateR ExampleSites.xlsx update:
#  1. Insert a new Sites row (ABT-162 / Cox Street bridge) between ABT-144 and ABT-237.
#  2. Add a new "Instructions" worksheet after "Sites" describing the Sites-tab fields.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sites sheet: insert the new monitoring location row (new row 6)
# ---------------------------------------------------------------------------
$sites = $wb.Worksheets.Item("Sites")

$sites.Rows.Item(6).Insert()
$sites.Range("A6").Value = "ABT-162"
$sites.Range("B6").Value = "Cox Street bridge"
$sites.Range("C6").Value = 42.399797
$sites.Range("D6").Value = -71.545985000000002
$sites.Range("E6").Value = "Assabet"

# restore the selection Excel leaves behind after this kind of edit
$sites.Range("B15").Select()

# ---------------------------------------------------------------------------
# 2. Add the "Instructions" worksheet right after "Sites"
# ---------------------------------------------------------------------------
$instr = $wb.Worksheets.Add($null, $sites)
$instr.Name = "Instructions"

# -- column widths ----------------------------------------------------------
$instr.Columns.Item(1).ColumnWidth = 30.29
$instr.Columns.Item(2).ColumnWidth = 96.43
$instr.Columns.Item(3).ColumnWidth = 21.43
$instr.Columns.Item(4).ColumnWidth = 21.57
$instr.Columns.Item(5).ColumnWidth = 25

# -- title / update-date rows ------------------------------------------------
$instr.Range("A1").Value = "The Sites tab must be formatted exactly like the Sites template, with all of the following fields."
$instr.Range("A2").Value = "The Sites tab must be the first tab in this workbook."
$instr.Range("C1").Value = "Template updated 5/19/23"
$instr.Range("C2").Value = "Samples updated 5/19/23"

$instr.Range("A1:A2").Font.Bold = $true
$instr.Range("C1:C2").Font.Color = 255

# -- header row (row 4) -------------------------------------------------------
$instr.Range("A4").Value = "Field"
$instr.Range("B4").Value = "Instructions"
$instr.Range("C4").Value = "Example"
$instr.Range("D4").Value = "Available Values"
$instr.Range("E4").Value = "Required?"

$hdr = $instr.Range("A4:E4")
$hdr.HorizontalAlignment = -4108
$hdr.Borders(7).LineStyle = 1
$hdr.Borders(7).Weight = 2
$hdr.Borders(10).LineStyle = 1
$hdr.Borders(10).Weight = 2
$hdr.Borders(8).LineStyle = 1
$hdr.Borders(8).Weight = 2
$hdr.Borders(9).LineStyle = 1
$hdr.Borders(9).Weight = -4138

# -- data rows 5-9 ------------------------------------------------------------
$fields = @("Monitoring Location ID", "Monitoring Location Name", "Monitoring Location Latitude ", "Monitoring Location Longitude", "Location Group")
$instructions = @(
  "Location ID that is used in your Results file.  Must match exactly.",
  "Name of monitoring location.",
  "Latitude of monitoring location in decimal form.  At least 5 decimals.",
  "Longitude of monitoring location in decimal form.  At least 5 decimals.",
  "An optional free-form grouping attribute.  This will allow you to summarize locations by group in the graphing and mapping analysis functions."
)
$required = @("Required", "Required for WQX", "Required for mapping", "Required for mapping", "Optional")

for ($i = 0; $i -lt 5; $i++) {
  $r = 5 + $i
  $instr.Range("A$r").Value = $fields[$i]
  $instr.Range("B$r").Value = $instructions[$i]
  $instr.Range("D$r").Value = "any"
  $instr.Range("E$r").Value = $required[$i]
}
$instr.Range("C5").Value = "ABT-010"
$instr.Range("C6").Value = "477 Lowell Rd, Concord"
$instr.Range("C7").Value = 42.470370000000003
$instr.Range("C8").Value = -71.362578999999997
$instr.Range("C9").Value = "Lower Assabet"

$body = $instr.Range("A5:E9")
$body.Borders(7).LineStyle = 1
$body.Borders(7).Weight = 2
$body.Borders(10).LineStyle = 1
$body.Borders(10).Weight = 2
$body.Borders(8).LineStyle = 1
$body.Borders(8).Weight = 2
$body.Borders(9).LineStyle = 1
$body.Borders(9).Weight = 2

# row 5's Field cell (A5) sits right under the header's thick rule, so it
# has no top edge of its own
$instr.Range("A5").Borders(8).LineStyle = -4142

$instr.Range("A5:E9").VerticalAlignment = -4160
$instr.Range("B5:B9").WrapText = $true
$instr.Range("C5:C9").HorizontalAlignment = -4108
$instr.Range("D5:E9").HorizontalAlignment = -4108
$instr.Range("D5:E9").Font.Italic = $true

$instr.Rows.Item(9).RowHeight = 30

# -- trailing bold empty cell -------------------------------------------------
$instr.Range("B11").Font.Bold = $true

# -- frozen header pane / view state ------------------------------------------
$instr.Range("B5").Select()
$excel.ActiveWindow.FreezePanes = $false
$instr.Range("B5").Select()
$excel.ActiveWindow.SplitColumn = 1
$excel.ActiveWindow.SplitRow = 4
$excel.ActiveWindow.FreezePanes = $true
$instr.Range("C3").Select()

Write-Output "done"
